$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two "report release date" header labels (shared strings)
$ws.Range("I9").Value = "1402-04-06 (11)"
$ws.Range("M9").Value = "1402-04-06 (3)"

# Update the quarterly figures in column M (newest quarter column)
$ws.Range("M14").Value = -8109
$ws.Range("M17").Value = 6093
$ws.Range("M19").Value = 44715
$ws.Range("M20").Value = 45340
$ws.Range("M21").Value = 385
$ws.Range("M22").Value = 45726
$ws.Range("M24").Value = 45726
